$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Insert 4 new rows before row 73 (shifts existing rows 73+ down by 4,
# and copies the formatting of row 72 into the new rows 73-76).
$ws.Rows("73:76").Insert()

# Fill in the new LOOKUP instruction description lines in column U,
# rows 72-75 (row 71 already holds the existing STORE note, row 76
# stays blank as a spacer row, matching the original layout).
# Order of entry matches the original authoring order so that new
# shared-string ids come out the same as in the target workbook.
$ws.Range("U72").Value2 = "LOOKUP: Searches for sub-sequence pointed by R[N.SRC], length R[N.SRC+1] in sequence pointed by R[N.RS], length R[N.SRC+1],"
$ws.Range("U74").Value2 = "if found, store its pointer to R[N.DST] and 0 to R[N.DST+1]"
$ws.Range("U75").Value2 = "if not found, store 1 to R[N.DST+1]"
$ws.Range("U73").Value2 = "ST_RS and ST_SRC are used to advance pointers of main and sub-sequences in loop for comparison."

# Merge C71:R75 into the big note cell and style it (left/top aligned,
# wrapped text) to hold the long "ALU vector mode" note.
$note = $ws.Range("C71:R75")
$note.Merge()
$note.HorizontalAlignment = -4131
$note.VerticalAlignment = -4160
$note.WrapText = $true
$note.Value2 = "Note: basic ALU operations have scalar mode (VL == 0) or vector mode (VL > 0). In vector mode ALU performs VL operations on its arguments. Argument will be read from register, if its stride is 0, or from memory, next argument value will be fetched from next register or memory, incrementing pointer by its stride. Same for destination."

# Update the selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("C76").Select()
